$wb = $excel.ActiveWorkbook

# --- Shared string update: "Ready for handoff" -> "In Translation" ---
# This status string is used on the Overview sheet (columns E/F, the
# per-language status columns) and on each language sheet's "Status"
# column (column C). Updating every cell that held the old text keeps
# a single shared-string entry for the new text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width updates (status columns narrower after text change) ---
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
